$wb = $excel.ActiveWorkbook

# Citywide Totals (sheet1.xml) - 21 changes
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('G2').Value = 91
$ws.Range('J2').Value = 129
$ws.Range('H3').Value = 161
$ws.Range('I3').Value = 198
$ws.Range('B6').Value = 391
$ws.Range('C6').Value = 506
$ws.Range('D6').Value = 437
$ws.Range('E6').Value = 508
$ws.Range('F6').Value = 574
$ws.Range('G6').Value = 446
$ws.Range('H6').Value = 476
$ws.Range('J6').Value = 432
$ws.Range('B7').Value = 530
$ws.Range('C7').Value = 666
$ws.Range('D7').Value = 680
$ws.Range('E7').Value = 746
$ws.Range('F7').Value = 826
$ws.Range('G7').Value = 691
$ws.Range('H7').Value = 765
$ws.Range('I7').Value = 862
$ws.Range('J7').Value = 826

# By Neighborhood (sheet2.xml) - 32 changes
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('B8').Value = 32
$ws.Range('C20').Value = 5
$ws.Range('B28').Value = 39
$ws.Range('C32').Value = 40
$ws.Range('B35').Value = 9
$ws.Range('C36').Value = 41
$ws.Range('H36').Value = 39
$ws.Range('J43').Value = 6
$ws.Range('E48').Value = 7
$ws.Range('J49').Value = 5
$ws.Range('F50').Value = 26
$ws.Range('D53').Value = 78
$ws.Range('I53').Value = 128
$ws.Range('J53').Value = 128
$ws.Range('H67').Value = 7
$ws.Range('H68').Value = 3
$ws.Range('G70').Value = 15
$ws.Range('D72').Value = 6
$ws.Range('J72').Value = 5
$ws.Range('G76').Value = 19
$ws.Range('C85').Value = 17
$ws.Range('H85').Value = 5
$ws.Range('E91').Value = 8
$ws.Range('B98').Value = 530
$ws.Range('C98').Value = 666
$ws.Range('D98').Value = 680
$ws.Range('E98').Value = 746
$ws.Range('F98').Value = 826
$ws.Range('G98').Value = 691
$ws.Range('H98').Value = 765
$ws.Range('I98').Value = 862
$ws.Range('J98').Value = 826

# Rogers Park (sheet3.xml) - 2 changes
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('G2').Value = 4
$ws.Range('G7').Value = 19

# Austin (sheet7.xml) - 2 changes
$ws = $wb.Worksheets.Item('Austin')
$ws.Range('B5').Value = 23
$ws.Range('B6').Value = 32

# Chicago Lawn (sheet9.xml) - 2 changes
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('C5').Value = 4
$ws.Range('C6').Value = 5

# Garfield Park (sheet10.xml) - 2 changes
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('C6').Value = 35
$ws.Range('C7').Value = 40

# Grand Crossing (sheet11.xml) - 4 changes
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('C6').Value = 36
$ws.Range('H6').Value = 25
$ws.Range('C7').Value = 41
$ws.Range('H7').Value = 39

# Little Italy, UIC (sheet13.xml) - 2 changes
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('F5').Value = 24
$ws.Range('F6').Value = 26

# Englewood (sheet18.xml) - 2 changes
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('B6').Value = 33
$ws.Range('B7').Value = 39

# Loop (sheet22.xml) - 6 changes
$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J2').Value = 21
$ws.Range('I3').Value = 30
$ws.Range('D6').Value = 48
$ws.Range('D7').Value = 78
$ws.Range('I7').Value = 128
$ws.Range('J7').Value = 128

# West Loop (sheet24.xml) - 2 changes
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('E6').Value = 7
$ws.Range('E7').Value = 8

# Grand Boulevard (sheet25.xml) - 2 changes
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('B5').Value = 7
$ws.Range('B6').Value = 9

# Printers Row (sheet39.xml) - 4 changes
$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range('D4').Value = 6
$ws.Range('J4').Value = 3
$ws.Range('D5').Value = 6
$ws.Range('J5').Value = 5

# United Center (sheet52.xml) - 4 changes
$ws = $wb.Worksheets.Item('United Center')
$ws.Range('C4').Value = 14
$ws.Range('H4').Value = 2
$ws.Range('C5').Value = 17
$ws.Range('H5').Value = 5

# Lincoln Park (sheet56.xml) - 2 changes
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('E5').Value = 6
$ws.Range('E6').Value = 7

# Lincoln Square (sheet59.xml) - 2 changes
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('J2').Value = 3
$ws.Range('J5').Value = 5

# Irving Park (sheet61.xml) - 2 changes
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('J5').Value = 2
$ws.Range('J6').Value = 6

# Old Town (sheet67.xml) - 2 changes
$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('G5').Value = 12
$ws.Range('G6').Value = 15

# O'Hare (sheet72.xml) - 2 changes
$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range('G3').Value = 2
$ws.Range('G6').Value = 3

# Norwood Park (sheet86.xml) - 2 changes
$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range('F5').Value = 6
$ws.Range('F6').Value = 7
